$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row (row 9) for the new mail entry ---
$logs = $wb.Worksheets.Item("Logs")
$logs.Range("A9").Value = "Ruilen van product"
$logs.Range("B9").Value = "mailmind.test@zohomail.eu"
$logs.Range("C9").Value = "Kan ik dit product ruilen voor een andere maat?"
$logs.Range("D9").Value = "Retour / Terugbetaling"
$logs.Range("F9").Value = "2025-06-23 18:20:14"
$logs.Range("G9").Value = "Nee"

# Extend the existing conditional formatting rules (Categorie / Beantwoord
# columns) so they keep covering the whole data range including the new row.
$catFmts = $logs.Range("D2:D8").FormatConditions
for ($i = 1; $i -le $catFmts.Count; $i++) {
    $catFmts.Item($i).ModifyAppliesToRange($logs.Range("D2:D9"))
}
$answeredFmts = $logs.Range("G2:G8").FormatConditions
for ($i = 1; $i -le $answeredFmts.Count; $i++) {
    $answeredFmts.Item($i).ModifyAppliesToRange($logs.Range("G2:G9"))
}

# --- Dashboard sheet: append new aggregated row (row 7) ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A7").Value = "Retour / Terugbetaling"
$dash.Range("B7").Value = 1

# --- Update the bar chart's category/value series to extend through the new row ---
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$7"
$series.Values = "='Dashboard'!`$B`$2:`$B`$7"
